$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 298, shifting existing rows 298:394 down to 299:395
$ws.Rows("298").Insert()

# Populate the new weekly record in row 298
$ws.Range("A298").Value = 4
$ws.Range("B298").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C298").Value = "Los Lagos"
$ws.Range("D298").Value = 44985
$ws.Range("E298").Value = 10
$ws.Range("F298").Value = 100112043
$ws.Range("G298").Value = "Pepino ensalada"
$ws.Range("H298").Value = "Sin especificar"
$ws.Range("I298").Value = "Primera"
$ws.Range("J298").Value = 400
$ws.Range("K298").Value = 12000
$ws.Range("L298").Value = 12000
$ws.Range("M298").Value = 12000
$ws.Range("N298").Value = "$/caja 60 unidades"
$ws.Range("O298").Value = "Región de Arica y Parinacota"
$ws.Range("P298").Value = 200
$ws.Range("Q298").Value = 60
$ws.Range("R298").Value = "Hortaliza"
